# "Updated sprint backlog and switched person numbers."
# On the "Sprint 2 (M2)" sheet, swap the assignees for two pairs of tasks,
# add a new data point, and update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 2 (M2)")

# Row 4 "Edit and commit Person 3 class": Bhavesh -> Stephen
$ws.Range("B4").Value = "Stephen"
# Row 5 "Edit and commit Person 4 class": Stephen -> Bhavesh
$ws.Range("B5").Value = "Bhavesh"

# New hours-logged value for row 4
$ws.Range("E4").Value = 0

# Row 9 "Create Build File 2": Bhavesh -> Stephen
$ws.Range("B9").Value = "Stephen"
# Row 10 "Create Build File 3": Stephen -> Bhavesh
$ws.Range("B10").Value = "Bhavesh"

# Update the sheet's active cell/selection to B11
$ws.Range("B11").Select()
